# Refresh the "cryptos" price/volume table on Sheet1 (columns D = Price,
# E = Volume(1h)) with the latest scraped values, as produced by the
# scheduled GitHub Actions job. Two rows (33/34) also swap coin identity
# (Celestia <-> WEMIXToken) because the source ranking reordered them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look like plain numbers (e.g. "317.20")
# must stay as literal text (matching the original text cells) instead of
# being auto-converted to a numeric value by Excel. Setting NumberFormat
# to Text ("@") before the assignment keeps the exact string, and then we
# restore the default "Normal" style so no stray style index is left on
# the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '43.109.35'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '2.553.13'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue "D5" '317.20'
Set-TextValue "D6" '97.14'
$ws.Range("E6").Value = '  +3.14%  '
Set-TextValue "D7" '0.574'
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("E8").Value = '  -0.02%  '
Set-TextValue "D9" '0.544'
$ws.Range("E9").Value = '  +3.41%  '
Set-TextValue "D10" '35.65'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  +0.79%  '
Set-TextValue "D12" '7.48'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("E13").Value = '  -4.68%  '
$ws.Range("D14").Value = '2.944.85'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = '2.603.88'
$ws.Range("E15").Value = '  +3.38%  '
Set-TextValue "D16" '15.01'
$ws.Range("E16").Value = '  -1.96%  '
Set-TextValue "D17" '0.845'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '43.128.62'
$ws.Range("E18").Value = '  +1.57%  '
Set-TextValue "D19" '6.86'
$ws.Range("E19").Value = '  +4.81%  '
Set-TextValue "D20" '12.57'
$ws.Range("E20").Value = '  -3.15%  '
$ws.Range("D21").Value = '0.0₃0963'
$ws.Range("E21").Value = '  +0.74%  '
Set-TextValue "D22" '70.01'
$ws.Range("E22").Value = '  +0.14%  '
Set-TextValue "D23" '254.16'
$ws.Range("E23").Value = '  +1.61%  '
Set-TextValue "D24" '2.94'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  +3.06%  '
Set-TextValue "D26" '26.71'
$ws.Range("E26").Value = '  +1.28%  '
$ws.Range("E27").Value = '  +0.28%  '
Set-TextValue "D28" '2.45'
$ws.Range("E28").Value = '  +2.86%  '
Set-TextValue "D29" '40.63'
$ws.Range("E29").Value = '  +5.90%  '
Set-TextValue "D30" '10.27'
$ws.Range("E30").Value = '  +1.49%  '
Set-TextValue "D31" '5.83'
$ws.Range("E31").Value = '  -1.11%  '
Set-TextValue "D32" '155.57'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D33" '2.71'
$ws.Range("E33").Value = '  +3.40%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D34" '19.18'
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E36").Value = '  +2.72%  '
Set-TextValue "D37" '3.31'
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("E38").Value = '  +1.65%  '
Set-TextValue "D39" '2.43'
$ws.Range("E39").Value = '  +5.23%  '
$ws.Range("E40").Value = '  +0.09%  '
Set-TextValue "D41" '22.08'
$ws.Range("E41").Value = '  -7.15%  '
Set-TextValue "D42" '3.85'
$ws.Range("E42").Value = '  +2.34%  '
Set-TextValue "D43" '0.0303'
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").Value = '1.981.04'
$ws.Range("E46").Value = '  -1.18%  '
Set-TextValue "D47" '84.61'
$ws.Range("E47").Value = '  +0.29%  '
Set-TextValue "D48" '9.03'
$ws.Range("E48").Value = '  +2.71%  '
$ws.Range("D49").Value = '2.800.59'
$ws.Range("E49").Value = '  +1.42%  '
Set-TextValue "D50" '104.40'
$ws.Range("E50").Value = '  +2.60%  '
Set-TextValue "D51" '74.06'
$ws.Range("E51").Value = '  +1.40%  '
